$wb = $excel.ActiveWorkbook

# Add the new worksheet and move it to the end (after "ODI Batting")
$ws = $wb.Worksheets.Add()
$ws.Name = "ODI Batting Extra"
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch the worksheet reference by name; the reference returned from
# Add() re-seats by index once the sheet order changes, so it now points
# at the wrong sheet after Move().
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Columns A, C, D, E, F hold text values (even the numeric-looking ones
# like "4483", "2", "0", "17.57%"), so force Text number format first or
# Excel will silently coerce them to numeric/percentage cells.
$ws.Range("A1:A6").NumberFormat = "@"
$ws.Range("C1:F6").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Row 2 - match 4483
$ws.Range("A2").Value = "4483"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "2"
$ws.Range("D2").Value = "0"
$ws.Range("E2").Value = "11.11%"
$ws.Range("F2").Value = "NO"

# Row 3 - match 4484
$ws.Range("A3").Value = "4484"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "0"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "NO"

# Row 4 - match 4564
$ws.Range("A4").Value = "4564"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "4"
$ws.Range("D4").Value = "0"
$ws.Range("E4").Value = "17.57%"
$ws.Range("F4").Value = "NO"

# Row 5 - match 4565 (player did not bat, all stats blank)
$ws.Range("A5").Value = "4565"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "NO"

# Row 6 - match 4567
$ws.Range("A6").Value = "4567"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "3"
$ws.Range("D6").Value = "1"
$ws.Range("E6").Value = "17.14%"
$ws.Range("F6").Value = "NO"

# Restore the default "Normal" style on the data rows (2-6) - the Text
# number format above was only needed transiently so literal numeric
# strings didn't get reinterpreted as numbers/percentages.
$ws.Range("A2:F6").Style = "Normal"

# Style the header row to match the header row style used on the other
# sheets in this workbook: bold font, centered/top-aligned, thin border.
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").HorizontalAlignment = -4108
$ws.Range("A1:F1").VerticalAlignment = -4160
$ws.Range("A1:F1").Borders.LineStyle = 1
